$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.951.18"
$ws.Range("E2").Value = "  +2.15%  "
$ws.Range("D3").Value = "1.651.64"
$ws.Range("E3").Value = "  +2.84%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "214.78"
$ws.Range("E5").Value = "  +1.30%  "
$ws.Range("E6").Value = "  +2.22%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +2.84%  "
$ws.Range("D9").Value = "0.0616"
$ws.Range("E9").Value = "  +1.73%  "
$ws.Range("D10").Value = "20.23"
$ws.Range("E10").Value = "  +5.21%  "
$ws.Range("D11").Value = "0.0879"
$ws.Range("E11").Value = "  +3.03%  "
$ws.Range("D12").Value = "1.885.19"
$ws.Range("E12").Value = "  +2.90%  "
$ws.Range("D13").Value = "1.656.32"
$ws.Range("E13").Value = "  +3.17%  "
$ws.Range("D14").Value = "4.08"
$ws.Range("E14").Value = "  +2.05%  "
$ws.Range("E15").Value = "  +2.71%  "
$ws.Range("D16").Value = "65.16"
$ws.Range("E16").Value = "  +2.88%  "
$ws.Range("D17").Value = "26.957.33"
$ws.Range("D18").Value = "235.77"
$ws.Range("E18").Value = "  +2.01%  "
$ws.Range("E19").Value = "  +1.31%  "
$ws.Range("E20").Value = "  +1.43%  "
$ws.Range("E22").Value = "  +3.24%  "
$ws.Range("E23").Value = "  +4.21%  "
$ws.Range("D24").Value = "2.23"
$ws.Range("E24").Value = "  +2.80%  "
$ws.Range("D25").Value = "145.33"
$ws.Range("E25").Value = "  -1.04%  "
$ws.Range("E26").Value = "  +2.22%  "
$ws.Range("E27").Value = "  +0.77%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E29").Value = "  +2.59%  "
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("E31").Value = "  +1.74%  "
$ws.Range("D32").Value = "1.555.17"
$ws.Range("E32").Value = "  +4.00%  "
$ws.Range("E33").Value = "  +2.91%  "
$ws.Range("D34").Value = "3.09"
$ws.Range("E34").Value = "  +5.21%  "
$ws.Range("E35").Value = "  +9.62%  "
$ws.Range("D36").Value = "2.41"
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("D37").Value = "0.586"
$ws.Range("E37").Value = "  +4.67%  "
$ws.Range("D38").Value = "0.894"
$ws.Range("E38").Value = "  +9.14%  "
$ws.Range("E39").Value = "  +2.73%  "
$ws.Range("D40").Value = "6.00"
$ws.Range("E40").Value = "  +3.64%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").Value = "66.07"
$ws.Range("E42").Value = "  +8.57%  "
$ws.Range("D43").Value = "2.23"
$ws.Range("E43").Value = "  +2.30%  "
$ws.Range("D44").Value = "1.792.23"
$ws.Range("E44").Value = "  +2.68%  "
$ws.Range("D45").Value = "0.775"
$ws.Range("E45").Value = "  +2.39%  "
$ws.Range("E46").Value = "  -0.80%  "
$ws.Range("D47").Value = "89.82"
$ws.Range("E47").Value = "  +0.27%  "
$ws.Range("E48").Value = "  +2.02%  "
$ws.Range("E49").Value = "  +3.16%  "
$ws.Range("E50").Value = "  +0.81%  "
$ws.Range("D51").Value = "7.62"
$ws.Range("E51").Value = "  +2.64%  "
